$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark sitting alone in the empty paragraph
#    right under the title (it gets relocated later in this script).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) "Anunciante" -> "Usuario" everywhere it appears as a capitalised whole
#    word (12 occurrences across the requirement table).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Anunciante", $true, $true, $false, $false, $false, $true, 1, $false, "Usuario", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) The one remaining lower-case mention ("...da pagina do anunciante")
#    becomes "...da pagina do Usuario" (capitalised).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("da página do anunciante", $true, $true, $false, $false, $false, $true, 1, $false, "da página do Usuario", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Typo fix: "logrado" -> "logado".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("já está logrado no site", $true, $false, $false, $false, $false, $true, 1, $false, "já está logado no site", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Drop the trailing space after "...banco de dados ".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("9- Sistema guarda novos dados no banco de dados ", $true, $false, $false, $false, $false, $true, 1, $false, "9- Sistema guarda novos dados no banco de dados", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Shorten the "desiste de atualizar" alternate-flow cell: merge the two
#    runs into a single sentence ending at "...nao confirma." and drop the
#    rest of the old text about the cancel button / redirect.
# ---------------------------------------------------------------------------
$oldFull = "6- Usuario desiste de atualizar anúncio e não confirma. Usuario clica em " + [char]0x2018 + "cancelar" + [char]0x2019 + " e é redirecionado para a página de listagem de anúncios"
$newFull = "6- Usuario desiste de atualizar anúncio e não confirma."
$d.Content.Find.Execute($oldFull, $true, $false, $false, $false, $false, $true, 1, $false, $newFull, 2) | Out-Null

# ---------------------------------------------------------------------------
# 7) Re-plant the "_GoBack" bookmark right after that sentence (between the
#    run and the paragraph mark). A collapsed Range sitting exactly on a
#    paragraph mark confuses Bookmarks.Add in this host, so we temporarily
#    append a sentinel character, bookmark just before it, then remove the
#    sentinel again -- the bookmark stays put, now correctly adjacent to the
#    real text.
# ---------------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*desiste de atualizar*") {
        $target = $p
        break
    }
}

$sentinelPos = $d.Range($target.Range.End - 1, $target.Range.End - 1)
$sentinelPos.InsertAfter("~")

$bmPos = $target.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$sentinelRange = $d.Range($bmPos, $bmPos + 1)
$sentinelRange.Text = ""

Write-Host "Edits applied"
